$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.388.25"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "3.104.13"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'243.16"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'615.53"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("D8").Value = "'0.384"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "3.105.68"
$ws.Range("E10").Value = "  +13.78%  "
$ws.Range("E11").Value = "  -3.99%  "
$ws.Range("D12").Value = "'0.204"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "'5.60"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "'34.45"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").Value = "91.468.49"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D18").Value = "3.112.97"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'14.72"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "'5.77"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "'444.62"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("E24").Value = "  -7.00%  "
$ws.Range("D25").Value = "'5.78"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").Value = "'88.42"
$ws.Range("E26").Value = "  -3.39%  "
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +26.13%  "
$ws.Range("D31").Value = "'0.230"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "'0.165"
$ws.Range("E32").Value = "  -8.74%  "
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("D34").Value = "'9.28"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "'0.986"
$ws.Range("D36").Value = "'7.65"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "'26.12"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("B38").Value = "MantraDAO"
$ws.Range("C38").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D38").Value = "'4.00"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "'1.94"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").Value = "'486.40"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "'0.434"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").Value = "'3.41"
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("D44").Value = "'22.20"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'159.39"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").Value = "'0.696"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").Value = "'44.09"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("E51").Value = "  -4.28%  "
